$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: booking_id=0, Myrthe Polfliet, AMS -> CDG, 2023-12-10, DELAYED, passport match=0 ---
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "Myrthe"
$ws.Cells.Item(2,3).Value = "Polfliet"
$ws.Cells.Item(2,4).Value = "AMS"
$ws.Cells.Item(2,5).Value = "CDG"
$ws.Cells.Item(2,6).Value = 45270
$ws.Cells.Item(2,7).Value = "DELAYED"
$ws.Cells.Item(2,8).Value = 0

# --- Row 3: booking_id=1, Myrthe Polfliet, CDG -> AMS, 2023-12-15, ON TIME , passport match=0 ---
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "Myrthe"
$ws.Cells.Item(3,3).Value = "Polfliet"
$ws.Cells.Item(3,4).Value = "CDG"
$ws.Cells.Item(3,5).Value = "AMS"
$ws.Cells.Item(3,6).Value = 45275
$ws.Cells.Item(3,7).Value = "ON TIME "
$ws.Cells.Item(3,8).Value = 0

# --- Row 4: booking_id=2, Zeno Koenigs, AMS -> FCO, 2023-12-17, DELAYED, passport match=1 ---
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "Zeno"
$ws.Cells.Item(4,3).Value = "Koenigs"
$ws.Cells.Item(4,4).Value = "AMS"
$ws.Cells.Item(4,5).Value = "FCO"
$ws.Cells.Item(4,6).Value = 45277
$ws.Cells.Item(4,7).Value = "DELAYED"
$ws.Cells.Item(4,8).Value = 1

# --- Row 5: booking_id=3, Myrthe Polfliet, AMS -> FCO, 2023-12-23, ON TIME , passport match=0 ---
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "Myrthe"
$ws.Cells.Item(5,3).Value = "Polfliet"
$ws.Cells.Item(5,4).Value = "AMS"
$ws.Cells.Item(5,5).Value = "FCO"
$ws.Cells.Item(5,6).Value = 45283
$ws.Cells.Item(5,7).Value = "ON TIME "
$ws.Cells.Item(5,8).Value = 0

# --- Row 6: booking_id=4, Myrthe Polfliet, AMS -> BCN, 2024-01-10, ON TIME , passport match=0 ---
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "Myrthe"
$ws.Cells.Item(6,3).Value = "Polfliet"
$ws.Cells.Item(6,4).Value = "AMS"
$ws.Cells.Item(6,5).Value = "BCN"
$ws.Cells.Item(6,6).Value = 45301
$ws.Cells.Item(6,7).Value = "ON TIME "
$ws.Cells.Item(6,8).Value = 0

# --- Row 7 (new): booking_id=5, Myrthe Polfliet, BCN -> AMS, 2024-01-29, ON TIME , passport match=0 ---
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "Myrthe"
$ws.Cells.Item(7,3).Value = "Polfliet"
$ws.Cells.Item(7,4).Value = "BCN"
$ws.Cells.Item(7,5).Value = "AMS"
$ws.Cells.Item(7,6).Value = 45320
$ws.Cells.Item(7,7).Value = "ON TIME "
$ws.Cells.Item(7,8).Value = 0

# Update selection to match the new active cell
$ws.Range("E6").Select()
